$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1672
$ws.Range("F3").Value = 9570
$ws.Range("F4").Value = 777
$ws.Range("F6").Value = 245
$ws.Range("F10").Value = 1439
$ws.Range("F11").Value = 578
$ws.Range("F13").Value = 1508
$ws.Range("F15").Value = 313
$ws.Range("F17").Value = 153
$ws.Range("F18").Value = 91
$ws.Range("F19").Value = 408
$ws.Range("F20").Value = 1115
$ws.Range("F22").Value = 25
$ws.Range("F23").Value = 4
$ws.Range("F24").Value = 52
$ws.Range("F25").Value = 291
$ws.Range("G25").Value = 60
$ws.Range("F26").Value = 7
$ws.Range("F27").Value = 269
$ws.Range("F29").Value = 608
$ws.Range("F30").Value = 648
$ws.Range("F33").Value = 182
$ws.Range("F35").Value = 14
$ws.Range("F36").Value = 189
$ws.Range("F37").Value = 336
$ws.Range("F38").Value = 512
$ws.Range("F39").Value = 317
$ws.Range("F40").Value = 634
$ws.Range("F41").Value = 515
$ws.Range("F42").Value = 749
$ws.Range("F43").Value = 331
$ws.Range("F44").Value = 283
$ws.Range("F48").Value = 67

$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 701
$ws.Range("F12").Value = 62
$ws.Range("F16").Value = 51
$ws.Range("F19").Value = 980
$ws.Range("F20").Value = 37
$ws.Range("F21").Value = 1069
$ws.Range("F22").Value = 295
$ws.Range("F23").Value = 657
$ws.Range("F24").Value = 32
$ws.Range("F33").Value = 174
$ws.Range("F39").Value = 31
$ws.Range("F40").Value = 30

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 779
$ws.Range("G5").Value = "不可售"
$ws.Range("F6").Value = 162
$ws.Range("F7").Value = 2419
$ws.Range("F8").Value = 3728
$ws.Range("F9").Value = 24
$ws.Range("F11").Value = 132
$ws.Range("F12").Value = 121

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1672
$ws.Range("F4").Value = 779
$ws.Range("F5").Value = 9570
$ws.Range("F6").Value = 3728
$ws.Range("F7").Value = 777
$ws.Range("F8").Value = 132
$ws.Range("F9").Value = 132
$ws.Range("F12").Value = 245
$ws.Range("F15").Value = 701
$ws.Range("F16").Value = 1439
$ws.Range("F17").Value = 578
$ws.Range("F18").Value = 121
$ws.Range("F19").Value = 121
$ws.Range("F20").Value = 1508
$ws.Range("F22").Value = 313
$ws.Range("F23").Value = 62
$ws.Range("F24").Value = 153
$ws.Range("F25").Value = 1115
$ws.Range("F29").Value = 291
$ws.Range("G29").Value = 60
$ws.Range("F30").Value = 37
$ws.Range("F31").Value = 7
$ws.Range("F32").Value = 270
$ws.Range("F33").Value = 1069
$ws.Range("F34").Value = 295
$ws.Range("F35").Value = 608
$ws.Range("F36").Value = 648
$ws.Range("F37").Value = 32
$ws.Range("F40").Value = 336
$ws.Range("F41").Value = 512
$ws.Range("F42").Value = 317
$ws.Range("F44").Value = 634
$ws.Range("F45").Value = 515
$ws.Range("F46").Value = 749
$ws.Range("F47").Value = 331
$ws.Range("F51").Value = 31
